# "Generate Report for Handback"
#
# For the zh-cn and de-de handback sheets: once a handback has actually
# happened, the report needs to show the resulting target/handback files
# and the real handback timestamp (instead of the 0001-01-01 placeholder),
# and the overall status label changes from "Ready for handoff" to
# "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

# The "Status" shared string is reused by the Overview sheet (B/C columns)
# and by the zh-cn / de-de sheets' B column, so updating it once here
# updates every place it is displayed.
$ws0 = $wb.Worksheets.Item("Overview")
$ws0.Range("B2").Value = "Handed back: in sync with en-US"

$sheetInfo = @(
  @{
    SheetName = "zh-cn"
    HandbackDateTime = "2016-03-09 16:37:06"
    Xlf1 = "196e36c0-f367-4e10-8da4-0c33f413e95d.2b7c70c3c3456f4a88c10450eaaf0121750ec224.zh-cn.xlf"
    Xlf2 = "77d0257a-5b9b-44e5-bfef-e654fa0e4069.a4e531d16eee642c77b7a42b524ad0dabcb01011.zh-cn.xlf"
    XlfBaseUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0900e9a878e9fe5699c4a504319de1d1a600cd98/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/high"
  },
  @{
    SheetName = "de-de"
    HandbackDateTime = "2016-03-09 16:37:22"
    Xlf1 = "196e36c0-f367-4e10-8da4-0c33f413e95d.2b7c70c3c3456f4a88c10450eaaf0121750ec224.de-de.xlf"
    Xlf2 = "77d0257a-5b9b-44e5-bfef-e654fa0e4069.a4e531d16eee642c77b7a42b524ad0dabcb01011.de-de.xlf"
    XlfBaseUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d3d376bbc4db1842abb38d8fce2e014907312fd3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/high"
  }
)

$mdBaseUrl = "https://github.com/OpenLocalizationTest/oltest/blob/97567c2643fa00c7b6f3cf5365b4115bacdc276b/e2e"
$md1 = "196e36c0-f367-4e10-8da4-0c33f413e95d.md"
$md2 = "77d0257a-5b9b-44e5-bfef-e654fa0e4069.md"

foreach ($info in $sheetInfo) {
  $ws = $wb.Worksheets.Item($info.SheetName)

  # Status text (same shared string as Overview, already updated above).
  $ws.Range("B2").Value = "Handed back: in sync with en-US"
  $ws.Range("B3").Value = "Handed back: in sync with en-US"

  # "Latest Target File" (E) / "Latest Handback File" (F) now get filled in
  # with the same files that were handed off, now handed back.
  $ws.Range("E2").Value = $md1
  $ws.Hyperlinks.Add($ws.Range("E2"), "$mdBaseUrl/$md1", "", "", $md1)
  $ws.Range("E2").Font.Underline = 2
  $ws.Range("E2").Font.Color = 15570276

  $ws.Range("F2").Value = $info.Xlf1
  $ws.Hyperlinks.Add($ws.Range("F2"), "$($info.XlfBaseUrl)/$($info.Xlf1)", "", "", $info.Xlf1)
  $ws.Range("F2").Font.Underline = 2
  $ws.Range("F2").Font.Color = 15570276

  $ws.Range("E3").Value = $md2
  $ws.Hyperlinks.Add($ws.Range("E3"), "$mdBaseUrl/$md2", "", "", $md2)
  $ws.Range("E3").Font.Underline = 2
  $ws.Range("E3").Font.Color = 15570276

  $ws.Range("F3").Value = $info.Xlf2
  $ws.Hyperlinks.Add($ws.Range("F3"), "$($info.XlfBaseUrl)/$($info.Xlf2)", "", "", $info.Xlf2)
  $ws.Range("F3").Font.Underline = 2
  $ws.Range("F3").Font.Color = 15570276

  # "Latest Handback DateTime" (G) moves from the 0001-01-01 placeholder to
  # the real handback timestamp for both rows.
  $ws.Range("G2").Value = $info.HandbackDateTime
  $ws.Range("G3").Value = $info.HandbackDateTime
}
